$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("HR EXPENSES")

# --- TRANSPORTATION (rows 7-13) and TOTAL (row 14) ---
$ws.Range("C7").Value = 2263.4
$ws.Range("D7").Value = 2331.3000000000002
$ws.Range("E7").Value = 2996.24
$ws.Range("C8").Value = 1520.9938000000002
$ws.Range("D8").Value = 1569.3812000000003
$ws.Range("E8").Value = 1991.3472000000002
$ws.Range("C9").Value = 189.87
$ws.Range("D9").Value = 257.08
$ws.Range("E9").Value = 302.02
$ws.Range("C10").Value = 355.08929999999998
$ws.Range("D10").Value = 325.77
$ws.Range("E10").Value = 0
$ws.Range("C11").Value = 116.31
$ws.Range("D11").Value = 150.19649999999999
$ws.Range("E11").Value = 197.65899999999999
$ws.Range("C12").Value = 375
$ws.Range("D12").Value = 413.79039999999998
$ws.Range("E12").Value = 391.16160000000002
$ws.Range("C13").Value = 1298
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 560
$ws.Range("C14").Value = 3855.2631000000001
$ws.Range("D14").Value = 2716.2181
$ws.Range("E14").Value = 3442.1878000000006
$ws.Range("F14").Value = 0

# --- LODGING & MEALS (rows 17-20) and TOTAL (row 21) ---
$ws.Range("C17").Value = 856.91
$ws.Range("D17").Value = 932.29
$ws.Range("E17").Value = 1655.7251999999999
$ws.Range("C18").Value = 55.43
$ws.Range("D18").Value = 76.048699999999997
$ws.Range("E18").Value = 94.95
$ws.Range("C19").Value = 278.46000000000004
$ws.Range("D19").Value = 1575.7896000000001
$ws.Range("E19").Value = 1162.8220000000001
$ws.Range("C20").Value = 1036.5079999999998
$ws.Range("D20").Value = 916.89999999999986
$ws.Range("E20").Value = 666.71679999999992
$ws.Range("C21").Value = 2227.308
$ws.Range("D21").Value = 3501.0282999999999
$ws.Range("E21").Value = 3580.2140000000004
$ws.Range("F21").Value = 0

# --- MISCELLANEOUS (rows 24-27) and TOTAL (row 28) ---
# These cells previously had General number format; typing currency values into
# them causes Excel to auto-apply a matching currency format (new style).
$ws.Range("C24:E27").NumberFormat = '"$"#,##0.00_);("$"#,##0.00)'
$ws.Range("C24").Value = 1225.643
$ws.Range("D24").Value = 1463.816
$ws.Range("E24").Value = 1586.3450000000003
$ws.Range("C25").Value = 1420.2705000000001
$ws.Range("D25").Value = 120.98740000000001
$ws.Range("E25").Value = 3403.7421999999997
$ws.Range("C26").Value = 863.72439999999995
$ws.Range("D26").Value = 804.06600000000003
$ws.Range("E26").Value = 1132.1944000000001
$ws.Range("C27").Value = 53.847999999999999
$ws.Range("D27").Value = 458.86
$ws.Range("E27").Value = 734.21599999999989
$ws.Range("C28").Value = 3563.4859000000001
$ws.Range("D28").Value = 2847.7294000000002
$ws.Range("E28").Value = 6856.4976000000006
$ws.Range("F28").Value = 0

# --- Update selection to match the final saved state ---
$ws.Range("D10").Select()
